$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'diari-2024-09-26 19:25:00-El nou centre per la gent gran costarà 3,5 milions d''euros'
$ws.Range("B2").Value = 45561.82228713979
$ws.Range("C2").Value = 'diari'
$ws.Range("E2").Value = 45561.80902777778
$ws.Range("F2").Value = 'parroquies'
$ws.Range("H2").Value = 'El nou centre per la gent gran costarà 3,5 milions d''euros'
$ws.Range("I2").Value = 'https://www.diariandorra.ad/parroquies/240926/les-obres-nou-espai-per-gent-gran-d-encamp-costaran-3-5-milions-d-euros_158863.html'

# Row 3
$ws.Range("A3").Value = 'diari-2024-09-26 18:55:00-La policia comissa 5.800 paquets de tabac valorats en més de 22.700 euros al Pas'
$ws.Range("B3").Value = 45561.82228713979
$ws.Range("C3").Value = 'diari'
$ws.Range("E3").Value = 45561.78819444445
$ws.Range("F3").Value = 'nacional'
$ws.Range("H3").Value = 'La policia comissa 5.800 paquets de tabac valorats en més de 22.700 euros al Pas'
$ws.Range("I3").Value = 'https://www.diariandorra.ad/nacional/240926/policia-comissa-5-800-paquets-tabac-valorats-mes-22-700-euros-pas_158862.html'

# Row 4
$ws.Range("A4").Value = 'diari-2024-09-26 18:01:00-El servei de salut mental del SAAS ha atès 3.000 adults durant el 2023'
$ws.Range("B4").Value = 45561.82228713979
$ws.Range("C4").Value = 'diari'
$ws.Range("E4").Value = 45561.75069444445
$ws.Range("F4").Value = 'nacional'
$ws.Range("H4").Value = 'El servei de salut mental del SAAS ha atès 3.000 adults durant el 2023'
$ws.Range("I4").Value = 'https://www.diariandorra.ad/nacional/240926/servei-salut-mental-saas-ates-3-000-adults-durant-aquest-any_158861.html'

# Row 5
$ws.Range("A5").Value = 'diari-2024-09-26 17:55:00-Retrets de la minoria per la baixa execució d''inversions a la capital'
$ws.Range("B5").Value = 45561.82228713979
$ws.Range("C5").Value = 'diari'
$ws.Range("E5").Value = 45561.74652777778
$ws.Range("F5").Value = 'parroquies'
$ws.Range("H5").Value = 'Retrets de la minoria per la baixa execució d''inversions a la capital'
$ws.Range("I5").Value = 'https://www.diariandorra.ad/parroquies/240926/retrets-minoria-per-baixa-execucio-d-inversions-capital_158860.html'

# Row 6
$ws.Range("A6").Value = 'diari-2024-09-26 17:46:00-Torna el ‘Bingo Art’ amb una selecció d’obres de l’artista Alejandra Pereyra'
$ws.Range("B6").Value = 45561.82228713979
$ws.Range("C6").Value = 'diari'
$ws.Range("E6").Value = 45561.74027777778
$ws.Range("F6").Value = 'parroquies'
$ws.Range("H6").Value = 'Torna el ‘Bingo Art’ amb una selecció d’obres de l’artista Alejandra Pereyra'
$ws.Range("I6").Value = 'https://www.diariandorra.ad/parroquies/240926/torna-bingo-art-amb-seleccio-d-obres-l-artista-alejandra-pereyra_158859.html'

# Row 7
$ws.Range("A7").Value = 'diari-2024-09-26 17:31:00-Govern i Ski Andorra signen un conveni en prevenció de riscos del personal de pistes'
$ws.Range("B7").Value = 45561.82228713979
$ws.Range("C7").Value = 'diari'
$ws.Range("E7").Value = 45561.72986111111
$ws.Range("F7").Value = 'nacional'
$ws.Range("H7").Value = 'Govern i Ski Andorra signen un conveni en prevenció de riscos del personal de pistes'
$ws.Range("I7").Value = 'https://www.diariandorra.ad/nacional/240926/govern-i-ski-andorra-signen-conveni-materia-prevencio-riscos-personal-pistes_158858.html'

# Row 8
$ws.Range("A8").Value = 'diari-2024-09-26 17:17:00-Conferència sobre ètica judicial i mitjans de comunicació per als membres del CSJ'
$ws.Range("B8").Value = 45561.82228713979
$ws.Range("C8").Value = 'diari'
$ws.Range("E8").Value = 45561.72013888889
$ws.Range("F8").Value = 'nacional'
$ws.Range("H8").Value = 'Conferència sobre ètica judicial i mitjans de comunicació per als membres del CSJ'
$ws.Range("I8").Value = 'https://www.diariandorra.ad/nacional/240926/conferencia-sobre-etica-judicial-i-mitjans-comunicacio-per-als-membres-csj_158857.html'

# Row 9
$ws.Range("A9").Value = 'diari-2024-09-26 17:04:00-Andorra Endavant acusa el Govern de "gestió ineficaç dels diners públics" pel cub led'
$ws.Range("B9").Value = 45561.82228713979
$ws.Range("C9").Value = 'diari'
$ws.Range("E9").Value = 45561.71111111111
$ws.Range("F9").Value = 'nacional'
$ws.Range("H9").Value = 'Andorra Endavant acusa el Govern de "gestió ineficaç dels diners públics" pel cub led'
$ws.Range("I9").Value = 'https://www.diariandorra.ad/nacional/240926/montaner-critica-despesa-exagerada-cub-led-per-informar-l-acord-d-associacio_158856.html'

# Row 10
$ws.Range("A10").Value = 'diari-2024-09-26 15:46:00-Sant Julià reivindicarà el sector ramader amb la primera festa de la transhumància'
$ws.Range("B10").Value = 45561.82228713979
$ws.Range("C10").Value = 'diari'
$ws.Range("E10").Value = 45561.65694444445
$ws.Range("F10").Value = 'parroquies'
$ws.Range("H10").Value = 'Sant Julià reivindicarà el sector ramader amb la primera festa de la transhumància'
$ws.Range("I10").Value = 'https://www.diariandorra.ad/parroquies/240926/sant-julia-reivindicara-sector-ramader-amb-primera-festa-transhumancia_158851.html'

# Row 11
$ws.Range("A11").Value = 'diari-2024-09-26 15:05:00-Ordino conclou amb èxit la segona edició del Seminari Reserves de la Biosfera'
$ws.Range("B11").Value = 45561.82228713979
$ws.Range("C11").Value = 'diari'
$ws.Range("E11").Value = 45561.62847222222
$ws.Range("F11").Value = 'parroquies'
$ws.Range("H11").Value = 'Ordino conclou amb èxit la segona edició del Seminari Reserves de la Biosfera'
$ws.Range("I11").Value = 'https://www.diariandorra.ad/parroquies/240926/ordino-conclou-amb-exit-segona-edicio-seminari-reserves-biosfera_158855.html'

# Row 12
$ws.Range("A12").Value = 'diari-2024-09-26 14:49:00-Sessió de comú d''Andorra la Vella'
$ws.Range("B12").Value = 45561.82228713979
$ws.Range("C12").Value = 'diari'
$ws.Range("E12").Value = 45561.61736111111
$ws.Range("F12").Value = 'diari-tv'
$ws.Range("H12").Value = 'Sessió de comú d''Andorra la Vella'
$ws.Range("I12").Value = 'https://www.diariandorra.ad/diari-tv/retransmissions/comu-andorra-vella/240926/sessio-comu-d-andorra-vella_158853.html'

# Row 13
$ws.Range("A13").Value = 'diari-2024-09-26 13:59:00-La propietat de la Borda Mateu trenca les negociacions amb l''FC Andorra per l''estadi'
$ws.Range("B13").Value = 45561.82228713979
$ws.Range("C13").Value = 'diari'
$ws.Range("E13").Value = 45561.58263888889
$ws.Range("F13").Value = 'esports'
$ws.Range("H13").Value = 'La propietat de la Borda Mateu trenca les negociacions amb l''FC Andorra per l''estadi'
$ws.Range("I13").Value = 'https://www.diariandorra.ad/esports/240926/trenquen-les-negocacions-l-fc-andorra-i-propietat-per-l-estadi-borda-mateu_158852.html'

# Row 14
$ws.Range("A14").Value = 'diari-2024-09-26 13:21:00-Atropellada una dona de 50 anys a La Massana'
$ws.Range("B14").Value = 45561.82228713979
$ws.Range("C14").Value = 'diari'
$ws.Range("E14").Value = 45561.55625
$ws.Range("F14").Value = 'nacional'
$ws.Range("H14").Value = 'Atropellada una dona de 50 anys a La Massana'
$ws.Range("I14").Value = 'https://www.diariandorra.ad/nacional/240926/atropellada-dona-50-anys-massana_158850.html'

# Row 15
$ws.Range("A15").Value = 'diari-2024-09-26 13:06:00-El grup caní de policia s''entrena amb la Guardia di Confine suïssa'
$ws.Range("B15").Value = 45561.82228713979
$ws.Range("C15").Value = 'diari'
$ws.Range("E15").Value = 45561.54583333333
$ws.Range("F15").Value = 'nacional'
$ws.Range("H15").Value = 'El grup caní de policia s''entrena amb la Guardia di Confine suïssa'
$ws.Range("I15").Value = 'https://www.diariandorra.ad/nacional/240926/grup-cani-policia-i-bombers-s-entrena-amb-guardia-di-confine-suissa_158849.html'

# Row 16
$ws.Range("A16").Value = 'diari-2024-09-26 12:56:00-El PS proposa esborrar l''historial clínic dels supervivents de càncer, VIH i hepatitis'
$ws.Range("B16").Value = 45561.82228713979
$ws.Range("C16").Value = 'diari'
$ws.Range("E16").Value = 45561.53888888889
$ws.Range("F16").Value = 'nacional'
$ws.Range("H16").Value = 'El PS proposa esborrar l''historial clínic dels supervivents de càncer, VIH i hepatitis'
$ws.Range("I16").Value = 'https://www.diariandorra.ad/nacional/240926/ps-proposa-esborrar-l-historial-clinic-dels-supervivents-cancer-vih-i-hepatitis_158848.html'

# Row 17
$ws.Range("A17").Value = 'diari-2024-09-26 12:29:00-El preu de les importacions creix un 6,8% al juny'
$ws.Range("B17").Value = 45561.82228713979
$ws.Range("C17").Value = 'diari'
$ws.Range("E17").Value = 45561.52013888889
$ws.Range("F17").Value = 'nacional'
$ws.Range("H17").Value = 'El preu de les importacions creix un 6,8% al juny'
$ws.Range("I17").Value = 'https://www.diariandorra.ad/nacional/240926/preu-les-importacions-creix-6-8_158840.html'

# Row 18
$ws.Range("A18").Value = 'diari-2024-09-26 12:23:00-Stop Violències engega la campanya ''F*cked'' per reivindicar l''avortament'
$ws.Range("B18").Value = 45561.82228713979
$ws.Range("C18").Value = 'diari'
$ws.Range("E18").Value = 45561.51597222222
$ws.Range("F18").Value = 'nacional'
$ws.Range("H18").Value = 'Stop Violències engega la campanya ''F*cked'' per reivindicar l''avortament'
$ws.Range("I18").Value = 'https://www.diariandorra.ad/nacional/240926/stop-violencies-engega-seva-campanya-f-cked-per-reivindicar-dret-les-dones-avortar_158845.html'

# Row 19
$ws.Range("A19").Value = 'diari-2024-09-26 10:55:00-FEDA augmenta un 4,6% la part fix de les tarifes de les xarxes de calor i fred'
$ws.Range("B19").Value = 45561.82228713979
$ws.Range("C19").Value = 'diari'
$ws.Range("E19").Value = 45561.45486111111
$ws.Range("F19").Value = 'nacional'
$ws.Range("H19").Value = 'FEDA augmenta un 4,6% la part fix de les tarifes de les xarxes de calor i fred'
$ws.Range("I19").Value = 'https://www.diariandorra.ad/nacional/240926/les-tarifes-les-xarxes-calor-feda-augmenten-4-6_158839.html'

# Row 20
$ws.Range("A20").Value = 'diari-2024-09-26 06:30:00-“Un dels meus allotjaments és l’únic Starlight del país”'
$ws.Range("B20").Value = 45561.82228713979
$ws.Range("C20").Value = 'diari'
$ws.Range("E20").Value = 45561.27083333334
$ws.Range("F20").Value = 'la-contra'
$ws.Range("H20").Value = '“Un dels meus allotjaments és l’únic Starlight del país”'
$ws.Range("I20").Value = 'https://www.diariandorra.ad/la-contra/240926/dels-meus-allotjaments-l-unic-starlight-pais_158825.html'

# Remove row 21 (article no longer present in updated scrape)
$ws.Rows.Item(21).Delete()
